$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sample Info sheet: add a new "batch" column (F) set to 1 for every
# existing sample row (rows 2-25), matching the style of column A.
# ------------------------------------------------------------------
$wsSample = $wb.Worksheets.Item("Sample Info")

$wsSample.Range("A1").Copy($wsSample.Range("F1"))
$wsSample.Range("F1").Value = "batch"

$wsSample.Range("A2:A25").Copy($wsSample.Range("F2:F25"))
$wsSample.Range("F2:F25").Value = 1

# ------------------------------------------------------------------
# Calibration sheet: add a new "batch" column (D) set to 1 for every
# existing calibration standard row (rows 2-46). Column D already
# carries the correct default style, so no format copy is required.
# ------------------------------------------------------------------
$wsCal = $wb.Worksheets.Item("Calibration")

$wsCal.Range("D1").Value = "batch"
$wsCal.Range("D2:D46").Value = 1

# ------------------------------------------------------------------
# Update selections / active sheet to match the saved workbook state.
# ------------------------------------------------------------------
[void]$wsSample.Range("H31").Select()

$wsParams = $wb.Worksheets.Item("Parameters")
[void]$wsParams.Range("C15").Select()

[void]$wsCal.Activate()
[void]$wsCal.Range("D2:D46").Select()
